$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G31").Value = "Acierto"
$ws.Range("H31").Value = 2

$ws.Range("G33").Value = "Acierto"
$ws.Range("H33").Value = 0.57

$ws.Range("G34").Value = "Acierto"
$ws.Range("H34").Value = 2
